$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.216.20"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.033.85"
$ws.Range("E3").Value = "  +3.78%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.11"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.58"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.399"
$ws.Range("E9").Value = "  +6.67%  "
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.31"
$ws.Range("E12").Value = "  +7.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.865"
$ws.Range("E13").Value = "  +4.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.332.46"
$ws.Range("E14").Value = "  +3.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.40"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.53"
$ws.Range("E16").Value = "  +5.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.033.09"
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.162.01"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.85"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.87"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.50"
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.57"
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.84"
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("E30").Value = "  +5.16%  "
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0668"
$ws.Range("E33").Value = "  +8.78%  "
$ws.Range("E34").Value = "  +2.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.51"
$ws.Range("E35").Value = "  +10.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.52"
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0983"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.20"
$ws.Range("E42").Value = "  +8.09%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.20"
$ws.Range("E43").Value = "  +2.76%  "
$ws.Range("E44").Value = "  +2.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.73"
$ws.Range("E45").Value = "  +4.52%  "
$ws.Range("E46").Value = "  +4.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.390.53"
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.54"
$ws.Range("E48").Value = "  +6.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.18"
$ws.Range("E49").Value = "  +19.63%  "
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.53"
$ws.Range("E51").Value = "  +1.64%  "
